$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B19").Value = "Latex akzente"
$ws.Range("B19").Select()
